# Implement rolling horizon simulation: refresh the forecasted demand /
# standard deviation figures (and the resulting average-demand summary
# statistic) that a rolling horizon pass would recompute.

$wb = $excel.ActiveWorkbook

# --- Productdata sheet: average demand summary value ---
$wsProduct = $wb.Worksheets.Item("Productdata")
$wsProduct.Range("H2").Value = 17.59221003605769

# --- ForecastedAverageDemand sheet: rolling-horizon forecast values ---
$wsAvgDemand = $wb.Worksheets.Item("ForecastedAverageDemand")
$wsAvgDemand.Range("B5").Value = 87
$wsAvgDemand.Range("B6").Value = 116
$wsAvgDemand.Range("B7").Value = 106
$wsAvgDemand.Range("B8").Value = 90
$wsAvgDemand.Range("B9").Value = 91
$wsAvgDemand.Range("B10").Value = 94
$wsAvgDemand.Range("B11").Value = 99
$wsAvgDemand.Range("B12").Value = 111
$wsAvgDemand.Range("B14").Value = 118

# --- ForcastedStandardDeviation sheet: rolling-horizon forecast values ---
$wsStdDev = $wb.Worksheets.Item("ForcastedStandardDeviation")
$wsStdDev.Range("B5").Value = 10.875
$wsStdDev.Range("B6").Value = 21.75
$wsStdDev.Range("B7").Value = 23.1875
$wsStdDev.Range("B8").Value = 21.09375
$wsStdDev.Range("B9").Value = 22.0390625
$wsStdDev.Range("B10").Value = 23.1328125
$wsStdDev.Range("B11").Value = 24.556640625
$wsStdDev.Range("B12").Value = 27.6416015625
$wsStdDev.Range("B14").Value = 29.47119140625
